# Auto-generated edit script: updates the cryptos list data (Price / Volume(1h) columns,
# plus a row swap for Kaspa/ImmutableX) to match the commit "Updated cryptos list" diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '73.072.45'
$ws.Range('E2').Value = '  +5.98%  '
$ws.Range('D3').Value = '2.660.88'
$ws.Range('E3').Value = '  +6.68%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '608.50'
$ws.Range('E5').Value = '  +2.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.60'
$ws.Range('E6').Value = '  +3.80%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +3.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.176'
$ws.Range('E9').Value = '  +15.76%  '
$ws.Range('D10').Value = '2.658.54'
$ws.Range('E10').Value = '  +6.63%  '
$ws.Range('E11').Value = '  +1.27%  '
$ws.Range('E12').Value = '  +5.55%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.13'
$ws.Range('E13').Value = '  +3.18%  '
$ws.Range('E14').Value = '  +11.54%  '
$ws.Range('D15').Value = '3.134.08'
$ws.Range('E16').Value = '  +5.78%  '
$ws.Range('D17').Value = '72.929.61'
$ws.Range('E17').Value = '  +6.12%  '
$ws.Range('D18').Value = '2.653.96'
$ws.Range('E18').Value = '  +6.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '386.29'
$ws.Range('E19').Value = '  +7.56%  '
$ws.Range('E20').Value = '  +7.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.95'
$ws.Range('E21').Value = '  +6.00%  '
$ws.Range('E22').Value = '  +5.73%  '
$ws.Range('E23').Value = '  +23.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.67'
$ws.Range('E24').Value = '  +5.59%  '
$ws.Range('E25').Value = '  +7.21%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.00'
$ws.Range('E27').Value = '  +12.46%  '
$ws.Range('D28').Value = '2.795.19'
$ws.Range('E28').Value = '  +6.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.40%  '
$ws.Range('D30').Value = '0.0₃0980'
$ws.Range('E30').Value = '  +12.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '540.32'
$ws.Range('E31').Value = '  +6.77%  '
$ws.Range('E33').Value = '  +11.77%  '
$ws.Range('E34').Value = '  +5.16%  '
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '162.53'
$ws.Range('E36').Value = '  +0.66%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.41'
$ws.Range('E37').Value = '  +4.78%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.42'
$ws.Range('E38').Value = '  +10.34%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.114'
$ws.Range('E39').Value = '  -3.38%  '
$ws.Range('E40').Value = '  +2.51%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.87'
$ws.Range('E41').Value = '  +10.73%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.68'
$ws.Range('E42').Value = '  +17.67%  '
$ws.Range('E43').Value = '  +8.81%  '
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('E45').Value = '  +6.58%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.82'
$ws.Range('E46').Value = '  +2.96%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '152.36'
$ws.Range('E47').Value = '  +2.48%  '
$ws.Range('E48').Value = '  +5.18%  '
$ws.Range('E49').Value = '  +7.60%  '
$ws.Range('E50').Value = '  +11.26%  '
$ws.Range('D51').Value = '0.0₆0269'
$ws.Range('E51').Value = '  +11.42%  '
